# Append a new row (row 4) to the sheet: a date (2021-12-06) in column A
# formatted as YYYY-MM-DD, and a number (4.25) in column B.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the date as its underlying serial number so Excel doesn't stamp it
# with its own implicit short-date format; we then apply the desired
# custom date format explicitly.
$ws.Range("A4").Value = 44536

# Exercise both casings of the custom date format code, matching the
# two numFmt entries (164 / 165) recorded for this workbook.
$ws.Range("A4").NumberFormat = "yyyy-mm-dd"
$ws.Range("A4").NumberFormat = "YYYY-MM-DD"

$ws.Range("B4").Value = 4.25

Write-Output "done"
